# Weekly_Report.xlsx update: add Week 3 / Week 4 entries for "Team" and
# "Tyler" sheets, and re-label Tyler's existing Week 0/1/2 headers.

$wb = $excel.ActiveWorkbook

# ---- Tyler sheet -----------------------------------------------------
$tyler = $wb.Worksheets.Item("Tyler")

# Re-label the existing week headers with the fuller titles.
$tyler.Range("A1").Value = "Week 0 - Lightning Talk"
$tyler.Range("A5").Value = "Week 1 - Architectural Review"
$tyler.Range("A10").Value = "Week 2 - Microarchitectural Review"

# Week 3 section.
$tyler.Range("A13").Value = "Week 3 - Basic Building Blocks Implementation"
$tyler.Range("A14").Value = "cordic.sv, LUT.sv, counter.sv, cordic_iteration.sv, CORDIC_tb.sv"

# Week 4 section.
$tyler.Range("A16").Value = "Week 4 - Processor/Accelerator Implementation"
$tyler.Range("A17").Value = "hier.sv, hier_tb.sv (incomplete)"
$tyler.Range("A18").Value = "angle_label_unit, angle_label_unit.sv"

$tyler.Select()
$tyler.Range("B22").Select()

# ---- Team sheet --------------------------------------------------------
$team = $wb.Worksheets.Item("Team")

$team.Range("A14").Value = "Week 3"
$team.Range("A15").Value = "Basic Building Blocks Implementation"
$team.Range("A18").Value = "Week 4"
$team.Range("A19").Value = "Processor/Accelerator Implementation"

# Team becomes the active sheet/tab again.
$team.Select()
$team.Range("C32").Select()
